# Updates the cryptos price list (D = Price, E = Volume(1h)) with freshly
# scraped values, and re-sorts several coin pairs whose rank order flipped
# between scrapes (rows 30/31, 34/35, 36/37, 38/39, 40/42/43/44/45/46/47/48).
#
# D-column prices are written via a temporary "@" (text) number format so
# numeric-looking strings (e.g. "1.01", "3.351.07") are stored as literal
# text instead of being auto-coerced into floating point numbers / losing
# the thousands-dot formatting used by the source site. ClearFormats()
# afterwards drops that temporary format so the cell keeps the workbook's
# original (unstyled) look.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.040.77"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -8.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.364.74"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -6.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +1.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "377.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -9.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "118.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.351.07"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -6.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -12.71%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.638"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -17.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.131"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -28.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000272"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -19.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.20"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -12.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.955.20"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.82"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -10.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.135"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.409.39"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.88"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -12.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.84"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.520.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -7.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.978"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -13.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.44"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -17.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.95"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.65"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -12.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.70"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -14.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "31.55"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -9.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.88"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -13.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.38"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -16.08%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.54"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -8.33%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.35"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -8.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.105"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -10.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.38"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -13.76%  "
$ws.Range("B34").Value = "Dai"
$ws.Range("C34").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.145"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -10.35%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "53.32"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.10%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "34.99"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -14.17%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0418"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -15.37%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.127"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -13.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.54"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.11%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.97"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +13.96%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.38"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.41%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -11.77%  "
$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.94"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -9.54%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.86"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.79%  "
$ws.Range("B47").Value = "PEPE"
$ws.Range("C47").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0562"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -27.34%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.61"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +10.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.83"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -11.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.53"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -17.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.264"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -15.57%  "
